$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The header label in A1 changed from "color" to "pal" (colors of the
# current palette).
$ws.Range("A1").Value = "pal"

# Reset the cached selection/active cell to A1 (was E9 previously).
$ws.Range("A1").Select()
